# Weekly fruit/vegetable price update: a new week's record (2 rows, Primera
# and Segunda quality) is inserted at the top of the data block (before the
# existing row 332), pushing the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 332.
$ws.Rows("332:333").Insert()

# Row 332 - "Primera" quality entry for the new date.
$ws.Range("A332").Value = 3
$ws.Range("B332").Value = "Femacal de La Calera"
$ws.Range("C332").Value = "Coquimbo"
$ws.Range("D332").Value = 44461
$ws.Range("D332").NumberFormat = $ws.Range("D334").NumberFormat
$ws.Range("E332").Value = 5
$ws.Range("F332").Value = 100112008
$ws.Range("G332").Value = "Coliflor"
$ws.Range("H332").Value = "Sin especificar"
$ws.Range("I332").Value = "Primera"
$ws.Range("J332").Value = 900
$ws.Range("K332").Value = 600
$ws.Range("L332").Value = 600
$ws.Range("M332").Value = 600
$ws.Range("N332").Value = "$/unidad"
$ws.Range("O332").Value = "Provincia de Quillota"
$ws.Range("P332").Value = 600
$ws.Range("Q332").Value = 1
$ws.Range("R332").Value = "Hortaliza"

# Row 333 - "Segunda" quality entry for the new date.
$ws.Range("A333").Value = 3
$ws.Range("B333").Value = "Femacal de La Calera"
$ws.Range("C333").Value = "Coquimbo"
$ws.Range("D333").Value = 44461
$ws.Range("D333").NumberFormat = $ws.Range("D334").NumberFormat
$ws.Range("E333").Value = 5
$ws.Range("F333").Value = 100112008
$ws.Range("G333").Value = "Coliflor"
$ws.Range("H333").Value = "Sin especificar"
$ws.Range("I333").Value = "Segunda"
$ws.Range("J333").Value = 850
$ws.Range("K333").Value = 500
$ws.Range("L333").Value = 500
$ws.Range("M333").Value = 500
$ws.Range("N333").Value = "$/unidad"
$ws.Range("O333").Value = "Provincia de Quillota"
$ws.Range("P333").Value = 500
$ws.Range("Q333").Value = 1
$ws.Range("R333").Value = "Hortaliza"
